# Minor text tweaks on the "emitBranch()" code-listing slides:
#  - Slide 16: "emit()  // leaves" -> "emit()   // leaves" (extra space added)
#  - Slide 18: three "-> emit(if (condition) ..." lines each lose one leading space

$p = $ppt.ActivePresentation

# ---- Slide 16 --------------------------------------------------------
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(2)          # "Content Placeholder 2"
$tr16 = $sh16.TextFrame.TextRange
$para16 = $tr16.Paragraphs(6)        # "    emit()  // leaves " + "boolean" + " value on top of stack"

# Select the whole first run (chars 1-22) and rewrite it in place so the
# run isn't split/fragmented.
$run16 = $para16.Characters(1, 22)
$run16.Text = "    emit()   // leaves "

# ---- Slide 18 ----------------------------------------------------------
$s18 = $p.Slides.Item(18)
$sh18 = $s18.Shapes.Item(2)          # "Content Placeholder 2"
$tr18 = $sh18.TextFrame.TextRange

# Paragraph 7: "        " + "Symbol.equals" + "    -> emit(if (condition) "BE $label""
$para7 = $tr18.Paragraphs(7)
$run7 = $para7.Characters(22, 38)
$run7.Text = '   -> emit(if (condition) "BE $label"'

# Paragraph 9: "        " + "Symbol.notEqual" + "  -> emit(if (condition) "BNE $label""
$para9 = $tr18.Paragraphs(9)
$run9 = $para9.Characters(24, 37)
$run9.Text = ' -> emit(if (condition) "BNE $label"'

# Paragraph 11: "        " + "Symbol.lessThan" + "  -> emit(if (condition) "BL $label""
$para11 = $tr18.Paragraphs(11)
$run11 = $para11.Characters(24, 36)
$run11.Text = ' -> emit(if (condition) "BL $label"'
